$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.963.42'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '1.954.31'
$ws.Range("E3").Value = '  -0.70%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.31%  '

$ws.Range("E6").Value = '  -0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4864'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2934'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.92%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07022'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '107.35'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.29%  '

$ws.Range("D12").Value = '1.954.47'
$ws.Range("E12").Value = '  -0.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07758'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.360'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6989'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '278.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.39%  '

$ws.Range("D17").Value = '30.973.70'
$ws.Range("E17").Value = '  -0.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007754'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.06%  '

$ws.Range("E19").Value = '  -0.37%  '

$ws.Range("D20").Value = '2.205.01'
$ws.Range("E20").Value = '  -0.82%  '

$ws.Range("E21").Value = '  -0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.480'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.485'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.741'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.166'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1046'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.398'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.33%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.602'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.83%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.563'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.392'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04882'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7517'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.79%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.163'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.87%  '

$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01993'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.67%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.679'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.525'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.05'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.106'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8948'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '109.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.84%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4432'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.56%  '

$ws.Range("E46").Value = '  -0.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.799'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '991.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1249'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.226'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.55%  '
